$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(134).Insert()

$ws.Range("A134").Value = 11
$ws.Range("B134").Value = "Vega Monumental Concepción"
$ws.Range("C134").Value = "Bíobío"
$ws.Range("D134").Value = 44589
$ws.Range("D134").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E134").Value = 8
$ws.Range("F134").Value = 100114014
$ws.Range("G134").Value = "Betarraga"
$ws.Range("H134").Value = "Sin especificar"
$ws.Range("I134").Value = "Primera"
$ws.Range("J134").Value = 800
$ws.Range("K134").Value = 600
$ws.Range("L134").Value = 650
$ws.Range("M134").Value = 619
$ws.Range("N134").Value = "`$/paquete 5 unidades"
$ws.Range("O134").Value = "Región Metropolitana"
$ws.Range("P134").Value = 124
$ws.Range("Q134").Value = 5
$ws.Range("R134").Value = "Hortaliza"
